$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-25 Sunday" "2026-01-26 Monday"

Replace-Text "27×21=567" "65×58=3770"
Replace-Text "65×30=1950" "70×81=5670"
Replace-Text "83×98=8134" "49×56=2744"
Replace-Text "30×46=1380" "44×31=1364"
Replace-Text "98×19=1862" "74×90=6660"

Replace-Text "18×14=252" "60×34=2040"
Replace-Text "13×29=377" "61×47=2867"
Replace-Text "18×93=1674" "26×64=1664"
Replace-Text "58×54=3132" "39×49=1911"
Replace-Text "87×52=4524" "75×66=4950"

Replace-Text "57×55=3135" "43×93=3999"
Replace-Text "88×99=8712" "61×24=1464"
Replace-Text "18×34=612" "13×39=507"
Replace-Text "90×72=6480" "71×79=5609"
Replace-Text "82×16=1312" "95×98=9310"

Replace-Text "31×44=1364" "53×74=3922"
Replace-Text "35×26=910" "50×29=1450"
Replace-Text "67×30=2010" "61×90=5490"
Replace-Text "22×29=638" "76×61=4636"
Replace-Text "92×15=1380" "83×84=6972"

Replace-Text "43×11=473" "52×49=2548"
Replace-Text "59×81=4779" "18×71=1278"
Replace-Text "46×52=2392" "87×41=3567"
Replace-Text "68×88=5984" "41×78=3198"
Replace-Text "63×26=1638" "45×65=2925"
